$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 282 for the newest price-report week (2022-02-18)
$ws.Rows.Item(282).Resize(4).Insert()

# Row 282
$ws.Range("A282").Value = 8
$ws.Range("B282").Value = "Terminal La Palmera de La Serena"
$ws.Range("C282").Value = "Coquimbo"
$ws.Range("D282").Value = 44610
$ws.Range("E282").Value = 4
$ws.Range("F282").Value = "Fruta"
$ws.Range("G282").Value = 100103
$ws.Range("H282").Value = "Frutos de hueso (carozo)"
$ws.Range("I282").Value = 100103006
$ws.Range("J282").Value = "Nectarín"
$ws.Range("K282").Value = "August Red"
$ws.Range("L282").Value = "Especial"
$ws.Range("M282").Value = 10
$ws.Range("N282").Value = 355000
$ws.Range("O282").Value = 360000
$ws.Range("P282").Value = 357500
$ws.Range("Q282").Value = "$/bins (420 kilos)"
$ws.Range("R282").Value = "Región Metropolitana"
$ws.Range("S282").Value = 851
$ws.Range("T282").Value = 420

# Row 283
$ws.Range("A283").Value = 8
$ws.Range("B283").Value = "Terminal La Palmera de La Serena"
$ws.Range("C283").Value = "Coquimbo"
$ws.Range("D283").Value = 44610
$ws.Range("E283").Value = 4
$ws.Range("F283").Value = "Fruta"
$ws.Range("G283").Value = 100103
$ws.Range("H283").Value = "Frutos de hueso (carozo)"
$ws.Range("I283").Value = 100103006
$ws.Range("J283").Value = "Nectarín"
$ws.Range("K283").Value = "August Red"
$ws.Range("L283").Value = "Primera"
$ws.Range("M283").Value = 16
$ws.Range("N283").Value = 335000
$ws.Range("O283").Value = 340000
$ws.Range("P283").Value = 337500
$ws.Range("Q283").Value = "$/bins (420 kilos)"
$ws.Range("R283").Value = "Región Metropolitana"
$ws.Range("S283").Value = 804
$ws.Range("T283").Value = 420

# Row 284
$ws.Range("A284").Value = 8
$ws.Range("B284").Value = "Terminal La Palmera de La Serena"
$ws.Range("C284").Value = "Coquimbo"
$ws.Range("D284").Value = 44610
$ws.Range("E284").Value = 4
$ws.Range("F284").Value = "Fruta"
$ws.Range("G284").Value = 100103
$ws.Range("H284").Value = "Frutos de hueso (carozo)"
$ws.Range("I284").Value = 100103006
$ws.Range("J284").Value = "Nectarín"
$ws.Range("K284").Value = "Venus"
$ws.Range("L284").Value = "Primera"
$ws.Range("M284").Value = 16
$ws.Range("N284").Value = 330000
$ws.Range("O284").Value = 335000
$ws.Range("P284").Value = 332500
$ws.Range("Q284").Value = "$/bins (420 kilos)"
$ws.Range("R284").Value = "Región de O'Higgins"
$ws.Range("S284").Value = 792
$ws.Range("T284").Value = 420

# Row 285
$ws.Range("A285").Value = 8
$ws.Range("B285").Value = "Terminal La Palmera de La Serena"
$ws.Range("C285").Value = "Coquimbo"
$ws.Range("D285").Value = 44610
$ws.Range("E285").Value = 4
$ws.Range("F285").Value = "Fruta"
$ws.Range("G285").Value = 100103
$ws.Range("H285").Value = "Frutos de hueso (carozo)"
$ws.Range("I285").Value = 100103006
$ws.Range("J285").Value = "Nectarín"
$ws.Range("K285").Value = "Venus"
$ws.Range("L285").Value = "Segunda"
$ws.Range("M285").Value = 10
$ws.Range("N285").Value = 300000
$ws.Range("O285").Value = 305000
$ws.Range("P285").Value = 302500
$ws.Range("Q285").Value = "$/bins (420 kilos)"
$ws.Range("R285").Value = "Región de O'Higgins"
$ws.Range("S285").Value = 720
$ws.Range("T285").Value = 420
